# Updated cryptos list on Sun Feb 25 16:35:45 UTC 2024 with GitHub Actions
# Refresh the "Price" (column D) and "Volume(1h)" (column E) figures for the
# cryptocurrency rows on the active sheet to match the latest scrape.
#
# Column D values are stored as plain text in the source sheet (mixed
# formats like "51.443.81" or "0.0860" are not valid numbers), so force the
# text number format before writing so Excel doesn't re-interpret them as
# numeric values and silently normalise things like trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Cells.Item(2, 4).Value = '51.443.81'
$ws.Range("D2").Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  +0.61%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Cells.Item(3, 4).Value = '3.040.77'
$ws.Range("D3").Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  +2.67%  '
$ws.Cells.Item(4, 5).Value = '  +0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '385.06'
$ws.Range("D5").Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  +1.37%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '102.36'
$ws.Range("D6").Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.02%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = '0.543'
$ws.Range("D7").Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -0.26%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = '0.583'
$ws.Range("D9").Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -0.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '36.78'
$ws.Range("D10").Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  +0.76%  '
$ws.Cells.Item(11, 5).Value = '  +0.24%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = '0.0860'
$ws.Range("D12").Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  +0.33%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = '3.527.86'
$ws.Range("D13").Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  +3.06%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '18.68'
$ws.Range("D14").Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  +2.29%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = '7.73'
$ws.Range("D15").Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -0.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '3.044.19'
$ws.Range("D16").Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  +3.11%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = '0.969'
$ws.Range("D17").Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -2.27%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = '10.52'
$ws.Range("D18").Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -4.75%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = '51.544.23'
$ws.Range("D19").Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  +0.78%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '3.14'
$ws.Range("D20").Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  -2.62%  '
$ws.Cells.Item(21, 5).Value = '  -0.93%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.0₃0964'
$ws.Range("D22").Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  +0.69%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Cells.Item(23, 4).Value = '70.23'
$ws.Range("D23").Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  +0.24%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '268.71'
$ws.Range("D24").Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.80%  '
$ws.Cells.Item(25, 5).Value = '  -1.16%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = '8.21'
$ws.Range("D26").Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  +5.59%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = '26.91'
$ws.Range("D27").Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  +3.90%  '
$ws.Cells.Item(28, 5).Value = '  +3.73%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '7.23'
$ws.Range("D29").Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -3.29%  '
$ws.Cells.Item(30, 5).Value = '  -0.01%  '
$ws.Cells.Item(31, 5).Value = '  -1.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = '10.25'
$ws.Range("D32").Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -0.23%  '
$ws.Cells.Item(33, 5).Value = '  +0.40%  '
$ws.Cells.Item(34, 5).Value = '  +2.56%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = '50.44'
$ws.Range("D35").Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -1.51%  '
$ws.Cells.Item(36, 5).Value = '  +2.19%  '
$ws.Cells.Item(37, 5).Value = '  -0.01%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '3.33'
$ws.Range("D38").Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  +2.79%  '
$ws.Cells.Item(39, 5).Value = '  +7.12%  '
$ws.Cells.Item(40, 5).Value = '  +2.64%  '
$ws.Cells.Item(41, 5).Value = '  +1.77%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = '128.19'
$ws.Range("D42").Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  +2.39%  '
$ws.Cells.Item(43, 5).Value = '  -0.74%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '2.53'
$ws.Range("D44").Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  +0.12%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '3.73'
$ws.Range("D45").Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  +4.92%  '
$ws.Cells.Item(46, 5).Value = '  +1.25%  '
$ws.Cells.Item(47, 5).Value = '  +4.15%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '2.07'
$ws.Range("D48").Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  +2.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = '2.027.26'
$ws.Range("D49").Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -0.17%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = '3.344.65'
$ws.Range("D50").Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  +2.78%  '
$ws.Cells.Item(51, 5).Value = '  +6.34%  '
